$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "937÷9=104, 1" "641÷3=213, 2"
Replace-Text "272÷5=54, 2" "642÷4=160, 2"
Replace-Text "598÷8=74, 6" "682÷8=85, 2"
Replace-Text "432÷6=72, 0" "391÷5=78, 1"
Replace-Text "342÷9=38, 0" "916÷3=305, 1"
Replace-Text "707÷4=176, 3" "772÷8=96, 4"
Replace-Text "253÷9=28, 1" "297÷8=37, 1"
Replace-Text "681÷7=97, 2" "136÷6=22, 4"
Replace-Text "845÷3=281, 2" "353÷2=176, 1"
Replace-Text "116÷6=19, 2" "564÷2=282, 0"
Replace-Text "378÷8=47, 2" "231÷2=115, 1"
Replace-Text "427÷3=142, 1" "370÷3=123, 1"
Replace-Text "646÷8=80, 6" "744÷7=106, 2"
Replace-Text "332÷9=36, 8" "885÷7=126, 3"
Replace-Text "448÷4=112, 0" "521÷8=65, 1"
Replace-Text "129÷7=18, 3" "948÷2=474, 0"
Replace-Text "374÷2=187, 0" "268÷8=33, 4"
Replace-Text "148÷3=49, 1" "845÷3=281, 2"
Replace-Text "100÷2=50, 0" "106÷7=15, 1"
Replace-Text "940÷5=188, 0" "246÷7=35, 1"
Replace-Text "241÷5=48, 1" "887÷5=177, 2"
Replace-Text "951÷5=190, 1" "716÷9=79, 5"
Replace-Text "454÷4=113, 2" "590÷3=196, 2"
Replace-Text "194÷7=27, 5" "478÷4=119, 2"
Replace-Text "239÷3=79, 2" "650÷7=92, 6"
